$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old GAMKL label cell at L1 (it is replaced by headers at K2/Q2 below)
$ws.Range("L1").ClearContents()

# --- New "K" column: one-hot-ish encoded arrays for each gel position (rows 3-12) ---
# Write order matters: it determines the order new strings are appended to the
# shared-string table, so write these in the same sequence the source workbook used.
$ws.Range("K4").Value = "[1. 1. 0. 0. 0. 0. 0. 1. 0. 0. 1. 0.]"
$ws.Range("K5").Value = "[1. 0. 1. 0. 0. 0. 0. 1. 0. 0. 0. 1.]"
$ws.Range("K6").Value = "[1. 0. 0. 1. 0. 0. 0. 0. 1. 0. 1. 0.]"
$ws.Range("K7").Value = "[1. 0. 0. 0. 1. 0. 0. 0. 1. 0. 0. 1.]"

# Header above the new K column, centered
$ws.Range("K2").Value = "G A M K L"
$ws.Range("K2").HorizontalAlignment = -4108

$ws.Range("K8").Value = "[1. 0. 0. 0. 0. 1. 0. 0. 0. 1. 1. 0.]"
$ws.Range("K9").Value = "[1. 0. 0. 0. 0. 0. 1. 0. 0. 1. 0. 1.]"
$ws.Range("K10").Value = "[1. 0. 0. 0. 0. 0. 0. 0. 0. 0. 1. 0.]"
$ws.Range("K11").Value = "[1. 0. 0. 0. 0. 0. 0. 0. 0. 0. 0. 1.]"
$ws.Range("K12").Value = "[1. 0. 0. 0. 0. 0. 0. 1. 0. 0. 0. 0.]"

# K3 reuses the existing all-zero array string
$ws.Range("K3").Value = "[0. 0. 0. 0. 0. 0. 0. 0. 0. 0. 0. 0.]"

# The GAMKL label that used to live at L1 now also appears at Q2
$ws.Range("Q2").Value = "GAMKL"

# --- The old K/L/M (Gel label / filename / array) columns shift right by three
# columns, to P/Q/R, and move down one row (rows 3-10 -> rows 4-11) ---
$ws.Range("P4").Value = "Gel_A1"
$ws.Range("Q4").Value = "7.jpg"
$ws.Range("R4").Value = "[0. 0. 0. 0. 0. 0. 0. 0. 0. 0. 0. 0.]"

$ws.Range("P5").Value = "Gel_A1"
$ws.Range("Q5").Value = "8.jpg"
$ws.Range("R5").Value = "[0. 0. 0. 0. 0. 0. 0. 0. 0. 0. 0. 0.]"

$ws.Range("P6").Value = "Gel_A1"
$ws.Range("Q6").Value = "3.jpg"
$ws.Range("R6").Value = "[0. 0. 0. 0. 0. 0. 0. 0. 0. 0. 0. 0.]"

$ws.Range("P7").Value = "Gel_A1"
$ws.Range("Q7").Value = "4.jpg"
$ws.Range("R7").Value = "[0. 0. 0. 0. 0. 0. 0. 0. 0. 0. 0. 0.]"

$ws.Range("P8").Value = "Gel_A1"
$ws.Range("Q8").Value = "5.jpg"
$ws.Range("R8").Value = "[1. 0. 0. 0. 0. 0. 0. 0. 1. 0. 0. 0.]"

$ws.Range("P9").Value = "Gel_A1"
$ws.Range("Q9").Value = "0.jpg"
$ws.Range("R9").Value = "[0. 0. 0. 0. 0. 0. 0. 0. 0. 0. 0. 0.]"

$ws.Range("P10").Value = "Gel_A1"
$ws.Range("Q10").Value = "1.jpg"
$ws.Range("R10").Value = "[0. 0. 0. 0. 0. 0. 0. 0. 0. 0. 0. 0.]"

$ws.Range("P11").Value = "Gel_A1"
$ws.Range("Q11").Value = "2.jpg"
$ws.Range("R11").Value = "[1. 0. 0. 0. 0. 0. 0. 0. 1. 0. 0. 0.]"

# Clear the now-vacated old K/L/M cells for rows 3-10 (K3 was overwritten above,
# not cleared, since it keeps a value)
$ws.Range("L3").ClearContents()
$ws.Range("M3").ClearContents()
$ws.Range("L4").ClearContents()
$ws.Range("M4").ClearContents()
$ws.Range("L5").ClearContents()
$ws.Range("M5").ClearContents()
$ws.Range("L6").ClearContents()
$ws.Range("M6").ClearContents()
$ws.Range("L7").ClearContents()
$ws.Range("M7").ClearContents()
$ws.Range("L8").ClearContents()
$ws.Range("M8").ClearContents()
$ws.Range("L9").ClearContents()
$ws.Range("M9").ClearContents()
$ws.Range("L10").ClearContents()
$ws.Range("M10").ClearContents()

# Widen the new K column and tidy up the view
$ws.Columns.Item(11).ColumnWidth = 28.109375
$ws.Range("K15").Select()
